# Correction du registre de risque
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fix the three corrected risk-mitigation texts (typo / capitalization fixes) ---
$ws.Range("F4").Value = "S'assurer de bien maintenir les révisions à jour, travaille collaboratif mis fréquemment à jour et dont les changements sont réversibles au moyen d'une synchronisation sur un serveur web protégé (GitHub)"
$ws.Range("F5").Value = "Faire des tests à partir de différentes plateformes et à partir de différents systèmes d'exploitation de manière à s'assurer l'homogénéité dans le fonctionnement de Matlab"
$ws.Range("F6").Value = "Tester le simulateur sur le plus de versions différentes de PSIM"

# --- Stray leftover values below the table (D11:D14), mirroring D7:G7 stacked in one column ---
$ws.Range("D11").Value = $ws.Range("D7").Value2
$ws.Range("D12").Value = $ws.Range("E7").Value2
$ws.Range("D12").NumberFormat = "0%"
$ws.Range("D13").Value = $ws.Range("F7").Value2
$ws.Range("D14").Value = $ws.Range("G7").Value2

# --- D7 loses its table border (becomes a center/wrap style with no border) ---
$ws.Range("D7").Borders.LineStyle = -4142

# --- View state: selection moves to F5, scroll position resets to top-left A1 ---
$ws.Application.Goto($ws.Range("A1"))
$ws.Range("F5").Select() | Out-Null
